$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.449.01"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "3.087.91"
$ws.Range("E3").Value = "  +1.06%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  -1.27%  "

$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("E11").Value = "  -0.94%  "

$ws.Range("D12").Value = "3.615.63"
$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("E13").Value = "  +2.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.01%  "

$ws.Range("E15").Value = "  -1.78%  "

$ws.Range("D16").Value = "57.537.51"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").Value = "3.090.20"
$ws.Range("E17").Value = "  +1.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.96%  "

$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "334.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.84%  "

$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.502"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.83%  "

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0908"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("E32").Value = "  -3.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.08%  "

$ws.Range("E35").Value = "  -3.63%  "

$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("E37").Value = "  -1.50%  "

$ws.Range("E38").Value = "  -0.51%  "

$ws.Range("D39").Value = "3.128.53"
$ws.Range("E39").Value = "  +1.17%  "

$ws.Range("E40").Value = "  -0.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.670"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").Value = "2.288.79"
$ws.Range("E44").Value = "  +3.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0257"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.04%  "

$ws.Range("E46").Value = "  -1.21%  "

$ws.Range("E47").Value = "  -1.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.54%  "

$ws.Range("E49").Value = "  -3.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "253.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.24%  "

$ws.Range("E51").Value = "  +0.93%  "
